# Adds the "plate2" data block (rows 15-23) and its "max" summary row (25)
# to the "raw" sheet, mirroring the existing "plate1" block in rows 1-11,
# then restores the post-edit selection/view state on both the "raw" and
# "analysis" sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw")

# --- New "plate2" data block -------------------------------------------------
# Row 15: header row (plate2 / 1..12)
$ws.Range("A15").Value = "plate2"
$rowVals = @(1,2,3,4,5,6,7,8,9,10,11,12)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(15, 2 + $i).Value = $rowVals[$i]
}

# Rows 16-23: sample rows A-H
$ws.Range("A16").Value = "A"
$rowVals = @(0.08502052383417724,0.07001804572208457,0.07734701128128743,0.07298555411072875,0.07616323638089344,0.08131284584662313,0.07497706400245259,0.06799636856760724,0.07255447244824101,0.08584821202393504,0.081620525404382,0.09304297297899441)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(16, 2 + $i).Value = $rowVals[$i]
}

$ws.Range("A17").Value = "B"
$rowVals = @(0.06580191817995007,0.06859171679132477,0.06958418947873228,0.07487040496323345,0.06639865403283096,0.06709971392633152,0.07264508476042246,0.07250074048118683,0.06898673505641305,0.07139148085568361,0.0740003464262554,0.08543127423097355)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(17, 2 + $i).Value = $rowVals[$i]
}

$ws.Range("A18").Value = "C"
$rowVals = @(0.07757230801670459,0.07497383476016667,0.08010461346825262,0.06932654828951733,0.08044855055953194,0.0733460601146142,0.07753446330252961,0.0694162880706654,0.06763035413405605,0.0743894203830333,0.06527120270575477,0.07731339202846522)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(18, 2 + $i).Value = $rowVals[$i]
}

$ws.Range("A19").Value = "D"
$rowVals = @(0.07184851221996236,0.06956807834050512,0.07474314098067006,0.06800022388584497,0.07142452661786071,0.06925927355238348,0.07147483675879324,0.07508301752250912,0.07574517245990661,0.0696587253354611,0.07905121374070709,0.08153093461836997)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(19, 2 + $i).Value = $rowVals[$i]
}

$ws.Range("A20").Value = "E"
$rowVals = @(0.06959202407428752,0.0736113514563401,0.07113256595410215,0.07137773922581647,0.06738665081696148,0.07500347052471064,0.0802427860224443,0.0732198688984113,0.07175268776421531,0.07194631940684972,0.07678940459196101,0.0742283979903088)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(20, 2 + $i).Value = $rowVals[$i]
}

$ws.Range("A21").Value = "F"
$rowVals = @(0.09290608377263132,0.10438887766531174,0.09425029138633101,0.0740263904168839,0.07634017049410469,0.07056359159343592,0.06719737414382734,0.07227093735795413,0.07420766427018541,0.07121544495587362,0.07301407017669866,0.07265536568556952)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(21, 2 + $i).Value = $rowVals[$i]
}

$ws.Range("A22").Value = "G"
$rowVals = @(0.07878257846867993,0.081430808084095,0.07203330653274476,0.07528011074053462,0.08050150759673227,0.07589444751765331,0.07387004805336991,0.07215073276847532,0.07622189412337967,0.08992059585240868,0.07566526164916784,0.07971400750824675)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(22, 2 + $i).Value = $rowVals[$i]
}

$ws.Range("A23").Value = "H"
$rowVals = @(0.07906628516354877,0.07700001286213164,0.0822686547032136,0.22722319391486348,0.27522216101230257,0.21467380024475677,0.09465549729681523,0.08198318065091183,0.08430083447225302,0.09163798499867085,0.09633838338134354,0.09750034923944509)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(23, 2 + $i).Value = $rowVals[$i]
}

# --- Row 25: "max" summary row, mirroring row 11's MAX formulas over the ----
#     new plate2 block, with the same highlighted "max" label style as A11.
$ws.Range("A25").Value = "max"
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M")
foreach ($col in $cols) {
    $ws.Range($col + "25").Formula = "=MAX(" + $col + "16:" + $col + "23)"
}

# --- Restore selection / view state -----------------------------------------
$ws.Range("B25").Select() | Out-Null

$analysis = $wb.Worksheets.Item("analysis")
$analysis.Activate() | Out-Null
$analysis.Range("O5").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollColumn = 3
} catch {
}
